$wb = $excel.ActiveWorkbook

# --- Statistics sheet ---
$ws1 = $wb.Worksheets.Item("Statistics")

$stats = @(
    ,("2024-08-22 18:10:47", 50.78206774004448, 5)
    ,("2024-08-22 18:10:49", 50.26977158884897, 10)
    ,("2024-08-22 18:10:51", 52.23679678239277, 14)
    ,("2024-08-22 18:10:53", 54.26354674871339, 16)
    ,("2024-08-22 18:10:55", 47.62661735562251, 18)
    ,("2024-08-22 18:10:57", 47.38531304538915, 19)
    ,("2024-08-22 18:10:59", 35.94129370198269, 16)
    ,("2024-08-22 18:11:01", 36.5052821850792, 17)
    ,("2024-08-22 18:11:03", 45.46172357303922, 15)
    ,("2024-08-22 18:11:05", 34.00479330213628, 15)
    ,("2024-08-22 18:11:07", 43.25765093338282, 15)
    ,("2024-08-22 18:11:09", 42.29515067463542, 16)
    ,("2024-08-22 18:11:11", 52.44730393024552, 16)
    ,("2024-08-22 18:11:13", 42.77169918724547, 16)
    ,("2024-08-22 18:11:15", 34.48240939036887, 18)
    ,("2024-08-22 18:11:17", 41.22569760117101, 17)
    ,("2024-08-22 18:11:19", 48.16745391638945, 15)
    ,("2024-08-22 18:11:21", 47.40192046541741, 14)
)

for ($i = 0; $i -lt $stats.Count; $i++) {
    $r = $i + 2
    $row = $stats[$i]
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
}

# --- Accidents sheet ---
$ws2 = $wb.Worksheets.Item("Accidents")

$accidents = @(
    ,("2024-08-22 18:10:54", "Truck and Car", "42.43 and 32.19", 1)
    ,("2024-08-22 18:10:56", "Car and Car", "58.21 and 40.46", 2)
    ,("2024-08-22 18:10:59", "Car and Car", "0.00 and 65.00", 3)
    ,("2024-08-22 18:10:59", "Car and Car", "0.00 and 53.67", 4)
    ,("2024-08-22 18:11:02", "Car and Car", "55.15 and 0.00", 4)
    ,("2024-08-22 18:11:03", "Car and Car", "0.00 and 40.12", 5)
    ,("2024-08-22 18:11:03", "Car and Car", "58.20 and 0.00", 5)
    ,("2024-08-22 18:11:05", "Car and Car", "62.12 and 0.00", 5)
    ,("2024-08-22 18:11:07", "Car and Truck", "52.43 and 7.56", 6)
    ,("2024-08-22 18:11:08", "Car and Truck", "47.16 and 0.00", 6)
    ,("2024-08-22 18:11:13", "Car and Car", "0.04 and 55.48", 7)
    ,("2024-08-22 18:11:14", "Car and Truck", "0.00 and 50.47", 8)
    ,("2024-08-22 18:11:16", "Car and Car", "0.00 and 58.61", 9)
    ,("2024-08-22 18:11:18", "Car and Car", "0.00 and 49.92", 9)
    ,("2024-08-22 18:11:18", "Car and Car", "61.00 and 0.00", 9)
)

for ($i = 0; $i -lt $accidents.Count; $i++) {
    $r = $i + 2
    $row = $accidents[$i]
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 4).Value = $row[3]
}